$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.842.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.398.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.81"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "67.710.37"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.73"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "328.46"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.72"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -4.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.54"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.97"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "414.44"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.32%  "
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.95"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.58"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.26"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.36"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("E42").Value = "  -7.54%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0910"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("E48").Value = "  -7.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("E51").Value = "  -0.88%  "
